# Update UK data for Panel F15 sr 49
# - Extend the F-column "wave" shared-formula fill down to F60 (F60 = F58+1)
# - Append a new row 61: country=uk, survey_round=49, panel=F, wave=15,
#   date_recieved=2021-03-04, spss_name=20-100590_PFW15_Final_IntUse_nodups,
#   r_name=<derived formula>, received_final=1
# - Update the selection to H61 to mirror where the author ended up editing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60 (F column): pick up the shared "wave" increment formula that
# previously stopped at F58 (F34:F58 -> F34:F60). Value stays 15.
$ws.Range("F60").Formula = "=F58+1"

# New row 61 data
$ws.Range("A61").Value = 3
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = "uk"
$ws.Range("D61").Value = 49
$ws.Range("E61").Value = "F"
$ws.Range("F61").Value = 15
$ws.Range("G61").Value = "2021-03-04"
$ws.Range("H61").Value = "20-100590_PFW15_Final_IntUse_nodups"
$ws.Range("I61").Formula = "=C61&""_""&""sr""&TEXT(D61,""00"")&""_""&YEAR(G61)&TEXT(G61,""MM"")&TEXT(G61,""DD"")&""_p""&E61&""_wv""&TEXT(F61,""00"")&"""""
$ws.Range("J61").Value = 1

# Leave the selection on H61, matching where editing finished
$ws.Range("H61").Select()
